# Update the "想去人数" (want-to-go count) figures that changed between
# the two most recent data pulls for the 南宁-漫展信息 workbook.
#
# Sheet "展览" (index 1) and sheet "全部类型" (index 4) both list the same
# events in rows 5 and 6, so both need the same update to column F:
#   F5: 3429 -> 3434
#   F6: 348  -> 351

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F5").Value = 3434
    $ws.Range("F6").Value = 351
}
